$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for price cells that would otherwise be
# auto-converted to numbers by Excel (losing trailing zeros, etc.)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '67.567.77'
$ws.Range("E2").Value = '  -3.09%  '
$ws.Range("D3").Value = '3.270.22'
$ws.Range("E3").Value = '  -5.65%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '592.45'
$ws.Range("E5").Value = '  -3.06%  '
$ws.Range("E6").Value = '  -9.97%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.262.06'
$ws.Range("E8").Value = '  -5.81%  '
$ws.Range("D9").Value = '0.545'
$ws.Range("E9").Value = '  -8.36%  '
$ws.Range("E10").Value = '  -10.65%  '
$ws.Range("D11").Value = '6.74'
$ws.Range("E11").Value = '  -5.49%  '
$ws.Range("E12").Value = '  -10.51%  '
$ws.Range("D13").Value = '38.48'
$ws.Range("E13").Value = '  -13.52%  '
$ws.Range("D14").Value = '0.0000246'
$ws.Range("E14").Value = '  -8.99%  '
$ws.Range("D15").Value = '3.791.31'
$ws.Range("E15").Value = '  -5.70%  '
$ws.Range("D16").Value = '67.546.40'
$ws.Range("E16").Value = '  -3.17%  '
$ws.Range("D17").Value = '3.272.09'
$ws.Range("E17").Value = '  -5.52%  '
$ws.Range("D19").Value = '533.11'
$ws.Range("E19").Value = '  -8.81%  '
$ws.Range("D20").Value = '7.14'
$ws.Range("E20").Value = '  -13.04%  '
$ws.Range("D21").Value = '14.98'
$ws.Range("E21").Value = '  -12.99%  '
$ws.Range("D22").Value = '0.758'
$ws.Range("E22").Value = '  -11.16%  '
$ws.Range("D23").Value = '7.90'
$ws.Range("E23").Value = '  -12.34%  '
$ws.Range("D24").Value = '85.54'
$ws.Range("E24").Value = '  -10.85%  '
$ws.Range("D25").Value = '13.58'
$ws.Range("E25").Value = '  -11.03%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -10.57%  '
$ws.Range("D28").Value = '8.10'
$ws.Range("E28").Value = '  -6.44%  '
$ws.Range("E29").Value = '  -11.99%  '
$ws.Range("D30").Value = '29.23'
$ws.Range("E30").Value = '  -11.48%  '
$ws.Range("E31").Value = '  -4.64%  '
$ws.Range("D32").Value = '1.17'
$ws.Range("E32").Value = '  -6.16%  '
$ws.Range("E33").Value = '  -15.81%  '
$ws.Range("D34").Value = '5.74'
$ws.Range("E34").Value = '  -12.89%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '515.81'
$ws.Range("E36").Value = '  -12.20%  '
$ws.Range("D37").Value = '0.0444'
$ws.Range("E37").Value = '  -7.38%  '
$ws.Range("D38").Value = '53.43'
$ws.Range("E38").Value = '  -5.08%  '
$ws.Range("E39").Value = '  -10.72%  '
$ws.Range("D40").Value = '8.98'
$ws.Range("E40").Value = '  -15.45%  '
$ws.Range("E41").Value = '  -9.83%  '
$ws.Range("D42").Value = '2.78'
$ws.Range("E42").Value = '  -12.47%  '
$ws.Range("D43").Value = '2.947.21'
$ws.Range("E43").Value = '  -9.29%  '
$ws.Range("E44").Value = '  -10.06%  '
$ws.Range("D45").Value = '0.0₃0589'
$ws.Range("E45").Value = '  -15.89%  '
$ws.Range("D46").Value = '2.19'
$ws.Range("E46").Value = '  -9.22%  '
$ws.Range("D47").Value = '26.78'
$ws.Range("E47").Value = '  -13.35%  '
$ws.Range("E49").Value = '  -16.61%  '
$ws.Range("E50").Value = '  -9.82%  '
$ws.Range("D51").Value = '123.85'
$ws.Range("E51").Value = '  -7.60%  '
